# Re-run of the benchmark: Array Length values shift down by one
# (e.g. 1001 -> 1000) and the measured Seconds per Insert timings are
# refreshed with the latest run's numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Array Length (col A), new Seconds per Insert (col B).
# A $null in the third slot means column B is unchanged for that row.
$data = @(
    ,@(2, 1000, "3E-06")
    ,@(3, 2000, "2.9E-05")
    ,@(4, 3000, "1.3E-05")
    ,@(5, 4000, "1.4999999999999999E-05")
    ,@(6, 5000, "1.7999999999999997E-05")
    ,@(7, 6000, "2E-05")
    ,@(8, 7000, "2.2E-05")
    ,@(9, 8000, "2.5E-05")
    ,@(10, 9000, "2.6E-05")
    ,@(11, 10000, "2.6E-05")
    ,@(12, 11000, "3.1E-05")
    ,@(13, 12000, "3.5000000000000004E-05")
    ,@(14, 13000, "3.9E-05")
    ,@(15, 14000, "4.4E-05")
    ,@(16, 15000, "5.5E-05")
    ,@(17, 16000, "5.1E-05")
    ,@(18, 17000, "5.5E-05")
    ,@(19, 18000, "5E-05")
    ,@(20, 19000, "5.9999999999999995E-05")
    ,@(21, 20000, "6.500000000000001E-05")
    ,@(22, 21000, "5.9E-05")
    ,@(23, 22000, "8.599999999999999E-05")
    ,@(24, 23000, "9.3E-05")
    ,@(25, 24000, "9.6E-05")
    ,@(26, 25000, "0.000131")
    ,@(27, 26000, "0.00010899999999999999")
    ,@(28, 27000, "0.000102")
    ,@(29, 28000, "9.900000000000001E-05")
    ,@(30, 29000, "0.000114")
    ,@(31, 30000, "9.7E-05")
    ,@(32, 31000, "9.5E-05")
    ,@(33, 32000, "0.00011999999999999999")
    ,@(34, 33000, "0.000116")
    ,@(35, 34000, "0.000107")
    ,@(36, 35000, "0.000106")
    ,@(37, 36000, "0.000114")
    ,@(38, 37000, "0.000112")
    ,@(39, 38000, "0.000111")
    ,@(40, 39000, "0.000125")
    ,@(41, 40000, "0.000116")
    ,@(42, 41000, "0.000115")
    ,@(43, 42000, "0.000122")
    ,@(44, 43000, "0.000122")
    ,@(45, 44000, "0.000114")
    ,@(46, 45000, "0.000124")
    ,@(47, 46000, "0.000125")
    ,@(48, 47000, "0.00013800000000000002")
    ,@(49, 48000, "0.00013800000000000002")
    ,@(50, 49000, "0.000135")
    ,@(51, 50000, $null)
    ,@(52, 51000, "0.00013800000000000002")
    ,@(53, 52000, "0.000136")
    ,@(54, 53000, "0.000158")
    ,@(55, 54000, "0.00016700000000000002")
    ,@(56, 55000, "0.00016")
    ,@(57, 56000, "0.000163")
    ,@(58, 57000, "0.000156")
    ,@(59, 58000, "0.00016700000000000002")
    ,@(60, 59000, "0.000154")
    ,@(61, 60000, "0.000156")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 2).Value = [double]$row[2]
    }
}

Write-Output "Updated $($data.Count) rows"
